$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Defect")

$ws2.Range("A7").Value = "Save Object"
$ws2.Range("B7").Value = "Can't save object if the bag has not been enabled."
$ws2.Range("C7").Value = "UI"
$ws2.Range("D7").Value = "Fish"
$ws2.Range("F7").Value = "10 Mar"

$ws2.Range("I12").Select()
